$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" timestamp updated
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 11:28:37"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 11:28:25"
$wsZhCn.Range("K2").Value = "2016-09-05 11:29:31"

# de-de sheet: "Correspond Handback DateTime" updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-05 11:29:51"
